$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Populate the "Sales Order" table (rows 3-7) first so that the
# --- shared-string table fills up with the table's strings before the
# --- title string, matching the target string order.

# Header row (row 3)
$ws.Range("A3").Value = "Order ID"
$ws.Range("B3").Value = "Customer Name"
$ws.Range("C3").Value = "Product"
$ws.Range("D3").Value = "Quantity"
$ws.Range("E3").Value = "Price"

# Data row 4 (John Doe / Laptop) - highlighted
$ws.Range("A4").Value = 1001
$ws.Range("B4").Value = "John Doe"
$ws.Range("C4").Value = "Laptop"
$ws.Range("D4").Value = 2
$ws.Range("E4").Value = 750

# Data row 5 (Jane Smith / Keyboard)
$ws.Range("A5").Value = 1002
$ws.Range("B5").Value = "Jane Smith"
$ws.Range("C5").Value = "Keyboard"
$ws.Range("D5").Value = 3
$ws.Range("E5").Value = 50

# Data row 6 (Bob Johnson / Monitor)
$ws.Range("A6").Value = 1003
$ws.Range("B6").Value = "Bob Johnson"
$ws.Range("C6").Value = "Monitor"
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 200

# Data row 7 (Alice Brown / Mouse)
$ws.Range("A7").Value = 1004
$ws.Range("B7").Value = "Alice Brown"
$ws.Range("C7").Value = "Mouse"
$ws.Range("D7").Value = 5
$ws.Range("E7").Value = 25

# Title (row 1) - written last so its shared string gets the final index
$ws.Range("A1").Value = "Sales Order"
$ws.Range("A1").Font.Bold = $true

# Header row formatting: bold font + thin box border
$ws.Range("A3:E3").Font.Bold = $true
$ws.Range("A3:E3").Borders.LineStyle = 1

# Highlighted data row: thin box border + yellow fill
$ws.Range("A4:E4").Borders.LineStyle = 1
$ws.Range("A4:E4").Interior.Color = 65535

# Remaining data rows: thin box border only
$ws.Range("A5:E5").Borders.LineStyle = 1
$ws.Range("A6:E6").Borders.LineStyle = 1
$ws.Range("A7:E7").Borders.LineStyle = 1

# Widen column B to fit customer names
$ws.Columns.Item(2).ColumnWidth = 19

# Leave the active selection on G7, matching the saved view state
$ws.Range("G7").Select() | Out-Null
